$wb = $excel.ActiveWorkbook

# Clear the routing numbers for accounts rows 2-16 (the "J" column) -
# this data now comes from the bank instead of being stored per-account.
$accounts = $wb.Worksheets.Item("accounts")
$accounts.Range("J2:J16").ClearContents()

# Make "accounts" the active sheet/tab (it was "bills" before).
$accounts.Activate()
$accounts.Range("J2:J16").Select()
$accounts.Range("J16").Activate()
